$p = $ppt.ActivePresentation

# Slide 10: "Loopback Measurement Mode – Round-trip Delay Measurement Mode"
#           -> "Loopback (Round-trip Delay) Measurement Mode"
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(3)
$tr10 = $sh10.TextFrame.TextRange.Paragraphs(8)
$tr10.Runs(1).Text = "Loopback (Round-trip Delay) Measurement Mode"

# Slide 3: Title "Requirements and Scope" -> "Requirements, Goals and Scope"
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(1)
$sh3.TextFrame.TextRange.Runs(1).Text = "Requirements, Goals and Scope"

# Slide 4: content placeholder updates
$s4 = $p.Slides.Item(4)
$sh4 = $s4.Shapes.Item(2)

# "STAMP Extensions for SR moved to " -> "STAMP Extensions for SR was moved to "
$para2 = $sh4.TextFrame.TextRange.Paragraphs(2)
$para2.Runs(1).Text = "STAMP Extensions for SR was moved to "

# "Replaced TWAMP Light with STAMP draft" -> "Replaced TWAMP Light draft with STAMP draft"
$para3 = $sh4.TextFrame.TextRange.Paragraphs(3)
$para3.Runs(1).Text = "Replaced TWAMP Light draft with STAMP draft"
